# Apply the "filled in scoresheet" edit to the Scoresheet_Richter_Michael sheet.
#
# Intent of the diff: row 7 ("Hcp" row) previously only carried the
# player's handicap in B7 (32.2) with the 18 per-hole adjustment cells
# (C7:K7, M7:U7) left blank. The committed workbook fills those cells in
# with the player's actual per-hole strokes (identical to the "Par" row,
# row 3) and rounds the handicap in B7 down to a whole number (32). All
# the dependent formulas in rows 9-13 (and the running totals in L7/V7/W7)
# recalculate automatically from that.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scoresheet_Richter_Michael")

# Handicap, rounded from 32.2 to 32.
$ws.Range("B7").Value = 32

# Per-hole strokes for holes 1-9 (C7:K7) -- same values as the Par row (row 3).
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = 3
$ws.Range("I7").Value = 3
$ws.Range("J7").Value = 5
$ws.Range("K7").Value = 4

# Per-hole strokes for holes 10-18 (M7:U7).
$ws.Range("M7").Value = 4
$ws.Range("N7").Value = 5
$ws.Range("O7").Value = 3
$ws.Range("P7").Value = 4
$ws.Range("Q7").Value = 3
$ws.Range("R7").Value = 4
$ws.Range("S7").Value = 4
$ws.Range("T7").Value = 4
$ws.Range("U7").Value = 3

# The sheet's selection moved to J27 in the saved file.
$ws.Range("J27").Select()
